$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.586.97'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '3.439.86'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.49'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.98'
$ws.Range("E6").Value = '  +8.20%  '
$ws.Range("D7").Value = '3.439.56'
$ws.Range("E7").Value = '  +2.02%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.74'
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.126'
$ws.Range("E11").Value = '  +2.82%  '
$ws.Range("E12").Value = '  +1.70%  '
$ws.Range("D13").Value = '4.027.68'
$ws.Range("E13").Value = '  +2.15%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.122'
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.69'
$ws.Range("E15").Value = '  +8.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000174'
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").Value = '3.435.89'
$ws.Range("E17").Value = '  +2.07%  '
$ws.Range("D18").Value = '61.665.46'
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.27'
$ws.Range("E19").Value = '  +8.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.31'
$ws.Range("E20").Value = '  +3.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.48'
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '387.71'
$ws.Range("E22").Value = '  +3.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.564'
$ws.Range("E23").Value = '  +2.81%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.12'
$ws.Range("E24").Value = '  +3.13%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.77'
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000124'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '3.583.46'
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.180'
$ws.Range("E29").Value = '  -0.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.72'
$ws.Range("E30").Value = '  +3.99%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.53'
$ws.Range("E32").Value = '  -10.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.19'
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("E34").Value = '  +1.78%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.12'
$ws.Range("E36").Value = '  +2.42%  '
$ws.Range("D37").Value = '3.468.92'
$ws.Range("E37").Value = '  +2.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.00'
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.16'
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("E40").Value = '  +0.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '166.57'
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0781'
$ws.Range("E42").Value = '  +3.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.99'
$ws.Range("E43").Value = '  +10.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.795'
$ws.Range("E44").Value = '  +2.94%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.50'
$ws.Range("E46").Value = '  +2.97%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.36'
$ws.Range("E47").Value = '  +1.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.73'
$ws.Range("E48").Value = '  +1.31%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("E49").Value = '  -2.46%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.571.84'
$ws.Range("E50").Value = '  +4.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.90'
$ws.Range("E51").Value = '  +1.78%  '
